# Updated calibration settings for the WTI riskDriverType PnL factor calibrations workbook.
# In particular: considering 5 days moving average (non-standardized by std),
# less batches and less episodes but longer.
#
# This updates the fitted parameters (mu, B, sig2, c, p, mu_0/1, B_0/1, sig2_0/1,
# omega, alpha, beta, gamma) and the abs_epsi_autocorr series on each of the
# AR / SETAR / GARCH / TARCH / AR_TARCH worksheets.

$wb = $excel.ActiveWorkbook

# --- AR sheet ---
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value2 = 0.0004844214392533607
$ws.Range("B3").Value2 = 0.7798658526816847
$ws.Range("B4").Value2 = 0.1280320194750863
$ws.Range("B5").Value2 = "[0.9999999999999998, 0.26003763856789847, 0.20602449448478205, 0.1597469691832025, 0.20825986776328745, 0.3866647177859979, 0.21727082136495096, 0.22149372192605693, 0.19520850351139343, 0.1644898527715792, 0.1594244972717661, 0.13889031921796882, 0.19804811655713567, 0.18451236965685267, 0.18617461346997563, 0.17362687256743825, 0.14815905112567532, 0.13262519010195487, 0.15237850217607135, 0.15381158690601096]"

# --- SETAR sheet ---
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B3").Value2 = 0.908130081300813
$ws.Range("B4").Value2 = -0.1330743634244479
$ws.Range("B5").Value2 = 0.7003805508813929
$ws.Range("B6").Value2 = 0.0893598703809307
$ws.Range("B7").Value2 = 0.1507729279382539
$ws.Range("B8").Value2 = 0.629820546979082
$ws.Range("B9").Value2 = 0.07210868033083918
$ws.Range("B10").Value2 = "[1.0, 0.2574247320295584, 0.25225877390737717, 0.21420852687443614, 0.2129895607788758, 0.25274910329965766, 0.18989586675965722, 0.17883538724352283, 0.16787643220697704, 0.1567448687619555, 0.1806829002034267, 0.1526544170395751, 0.18293572393069413, 0.16468814120584158, 0.16681676155815273, 0.1639826028816237, 0.1638623156680472, 0.1604368708196883, 0.17092566112809743, 0.14391953812945438]"

# --- GARCH sheet ---
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value2 = -0.0003018641867290171
$ws.Range("B3").Value2 = 0.001258914372815595
$ws.Range("B4").Value2 = 0.06199770195676332
$ws.Range("B5").Value2 = 0.929033465014958
$ws.Range("B6").Value2 = "[1.0, 0.008357227091285778, -0.033522811890720365, -0.07715060906364973, -0.07048078881447314, 0.2702364925151802, 0.004314955160213761, 0.04211372372469305, 9.402743801944046e-05, -0.03467394783274602, -0.04626321607775001, -0.06763779991936066, -0.002035940512775672, -0.008138706890245801, 0.01687747015950369, 0.018778925322410527, -0.021237562162973293, -0.04839882685446983, -0.06910436267687978, -0.021792967432692673]"

# --- TARCH sheet ---
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value2 = -0.001933675721538139
$ws.Range("B3").Value2 = 0.001269239019146223
$ws.Range("B4").Value2 = 0.0557718394638874
$ws.Range("B5").Value2 = 0.9278333736867835
$ws.Range("B7").Value2 = 0.01501718335146807
$ws.Range("B6").Value2 = "[1.0, 0.006994302699577481, -0.0343544569674512, -0.07784624053438803, -0.07116840212161608, 0.27013063393374553, 0.0037843186126448104, 0.04147781650032528, -0.0005101822586543604, -0.035030644416461816, -0.04683547178875497, -0.0676438631901248, -0.002343650114997198, -0.008690267521857349, 0.016148648209967458, 0.018186364013949737, -0.021518107935262678, -0.04820059659914482, -0.06932555832395212, -0.021857769567628412]"

# --- AR_TARCH sheet ---
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value2 = 0.002665281902989896
$ws.Range("B3").Value2 = 0.001010051145472083
$ws.Range("B4").Value2 = 0.04647750030145447
$ws.Range("B5").Value2 = 0.934875304760624
$ws.Range("B7").Value2 = 0.02061546307598787
$ws.Range("B9").Value2 = 0.7868578141333236
$ws.Range("B6").Value2 = "[0.9999999999999998, 0.026910216166365526, -0.01813501093009271, -0.047260640713221856, -0.03339288538954057, 0.2142768596116096, -0.019706406178560662, 0.023322009328276168, -0.0017835395578737025, -0.04504212849192158, -0.03314730519719008, -0.06126497828459132, 0.004449291431843918, -0.00045798308374967735, 0.0051632238047335174, 0.014767889759065029, -0.021590259164887958, -0.04075613801926588, -0.04076433306592119, -0.023396020458887647]"
